# Applies "Some Corrections on Form submission" changes to the offer letter.
$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# --- Simple, unambiguous text replacements (unique counts verified against the diff) ---
Replace-All "03-10-2025" "07-10-2025"
Replace-All "Radhe Shyam" "Swati Sharma"
Replace-All "8845679893" "9567856785"
Replace-All "Radheshyam@gmail.com" "swati.sharma@rigvedit.com"
Replace-All "Solution Architect" "Senior Account Manager"
Replace-All "08-10-2025" "15-10-2025"
Replace-All "17,00,000" "15,00,000"
Replace-All "Seventeen Lakh Rupees Only" "Fifteen Lakh Rupees Only"
Replace-All "5,44,000" "4,80,000"
Replace-All "45,333" "40,000"
Replace-All "2,72,000" "2,40,000"
Replace-All "22,667" "20,000"
Replace-All "54,400" "48,000"
Replace-All "4,533" "4,000"
Replace-All "65,280" "57,600"
Replace-All "5,440" "4,800"
Replace-All "3,80,598" "3,20,392"
Replace-All "31,716" "26,699"
Replace-All "14,99,678" "13,29,392"
Replace-All "1,24,973" "1,10,783"
Replace-All "78,930" "69,968"
Replace-All "6,578" "5,831"
Replace-All "15,78,608" "13,99,360"
Replace-All "1,31,551" "1,16,613"
Replace-All "26,112" "23,040"
Replace-All "2,176" "1,920"
Replace-All "30,000" "20,000"
Replace-All "1,21,392" "1,00,640"
Replace-All "10,116" "8,387"
Replace-All "1,41,667" "1,25,000"
Replace-All "67,780" "60,100"
Replace-All "5,640" "5,000"

# --- Ambiguous value: "2,500" appears twice (Insurance Coverage monthly value,
# and Profession Tax annual value). Only the Insurance Coverage row's monthly
# value (row index 20, column index 2) changes to "1,667".
$table = $d.Tables.Item(1)
$cell = $table.Cell(21, 3)
$cell.Range.Text = "1,667"
